$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F - give it the same bold/centered/bordered style
# used by the other header cells (A1:E1) before writing its text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Trening"

# Full replacement data set (rows 2-13), columns A-F:
#   A = Timestamp (date/time serial number)
#   B = Seconds
#   C = Velocity
#   D = Acceleration_SMA
#   E = Velocity_Bin
#   F = Trening
$rows = @(
    @(2, 45684.59258171296, 599, 12.85, 2.793306521006993, "10-15", "Duża Gra"),
    @(3, 45684.59327731482, 659.1, 13.8, 2.448943887438093, "10-15", "Duża Gra"),
    @(4, 45684.59392893519, 715.4, 14.08, 2.541927678244454, "10-15", "Duża Gra"),
    @(5, 45684.59257824074, 598.7, 9.300000000000001, 2.337233253887721, "5-10", "Duża Gra"),
    @(6, 45684.59327268518, 658.7, 9.619999999999999, 2.060739278793335, "5-10", "Duża Gra"),
    @(7, 45684.59392430555, 715, 9.98, 1.98479664325714, "5-10", "Duża Gra"),
    @(8, 45684.59989421297, 1230.8, 12.23, 2.927401304244994, "10-15", "Mała Gra"),
    @(9, 45684.60002384259, 1242, 13.71, 3.556057657514301, "10-15", "Mała Gra"),
    @(10, 45684.60255162037, 1460.4, 10.99, 2.936538662229266, "10-15", "Mała Gra"),
    @(11, 45684.59989189815, 1230.6, 9.81, 2.697639737810408, "5-10", "Mała Gra"),
    @(12, 45684.60002037037, 1241.7, 9.029999999999999, 2.703704352889744, "5-10", "Mała Gra"),
    @(13, 45684.60254930556, 1460.2, 9.02, 2.836862632206507, "5-10", "Mała Gra")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}

# Apply the custom date/time number format to the Timestamp column (A2:A13).
# First set it on a single cell twice (lowercase then uppercase) to
# reproduce the numFmt 164/165 pair seen in the target workbook (164
# defined but unused by any style, 165 actually applied). Then apply the
# final format to the rest of the column so every cell shares the same
# style index (no extra cellXfs entries get created).
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A3:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
